$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'288.02"
$ws.Range("E2").Value = "'-1.04%"
$ws.Range("G2").Value = "'9"
$ws.Range("D3").Value = "'31.06"
$ws.Range("E3").Value = "'1.28%"
$ws.Range("G3").Value = "'9"
$ws.Range("D4").Value = "'4.930"
$ws.Range("E4").Value = "'-0.43%"
$ws.Range("G4").Value = "'9"
$ws.Range("D5").Value = "'0.07335"
$ws.Range("E5").Value = "'1.64%"
$ws.Range("G5").Value = "'9"
$ws.Range("D6").Value = "'2.217"
$ws.Range("E6").Value = "'20.32%"
$ws.Range("G6").Value = "'9"
$ws.Range("D7").Value = "'7.717"
$ws.Range("E7").Value = "'0.44%"
$ws.Range("G7").Value = "'9"
$ws.Range("D8").Value = "'3.733"
$ws.Range("E8").Value = "'-0.92%"
$ws.Range("G8").Value = "'9"
$ws.Range("D9").Value = "'0.9031"
$ws.Range("E9").Value = "'0.60%"
$ws.Range("G9").Value = "'9"
$ws.Range("D10").Value = "'0.09169"
$ws.Range("E10").Value = "'19.16%"
$ws.Range("G10").Value = "'9"
$ws.Range("D11").Value = "'0.1685"
$ws.Range("E11").Value = "'1.84%"
$ws.Range("G11").Value = "'9"
$ws.Range("D12").Value = "'0.08147"
$ws.Range("E12").Value = "'1.80%"
$ws.Range("G12").Value = "'9"
$ws.Range("E13").Value = "'2.83%"
$ws.Range("G13").Value = "'9"
$ws.Range("D14").Value = "'0.09951"
$ws.Range("E14").Value = "'-0.60%"
$ws.Range("G14").Value = "'9"
$ws.Range("D15").Value = "'0.001495"
$ws.Range("E15").Value = "'-0.37%"
$ws.Range("G15").Value = "'9"
$ws.Range("D16").Value = "'0.005745"
$ws.Range("E16").Value = "'0.65%"
$ws.Range("G16").Value = "'9"
$ws.Range("D17").Value = "'3.521"
$ws.Range("E17").Value = "'1.64%"
$ws.Range("G17").Value = "'9"
$ws.Range("D18").Value = "'2.076"
$ws.Range("E18").Value = "'-0.31%"
$ws.Range("G18").Value = "'9"
$ws.Range("D19").Value = "'0.3331"
$ws.Range("E19").Value = "'0.45%"
$ws.Range("G19").Value = "'9"
$ws.Range("E20").Value = "'-0.46%"
$ws.Range("G20").Value = "'9"
$ws.Range("D21").Value = "'4.184"
$ws.Range("E21").Value = "'3.53%"
$ws.Range("G21").Value = "'9"
$ws.Range("E22").Value = "'-12.03%"
$ws.Range("G22").Value = "'9"
$ws.Range("D23").Value = "'0.04536"
$ws.Range("E23").Value = "'0.71%"
$ws.Range("G23").Value = "'9"
$ws.Range("D24").Value = "'0.001209"
$ws.Range("E24").Value = "'-0.59%"
$ws.Range("G24").Value = "'9"
$ws.Range("G25").Value = "'9"
$ws.Range("E26").Value = "'3.91%"
$ws.Range("G26").Value = "'9"
$ws.Range("G27").Value = "'9"
$ws.Range("G28").Value = "'9"
$ws.Range("G29").Value = "'9"
$ws.Range("G30").Value = "'9"
$ws.Range("G31").Value = "'9"
$ws.Range("G32").Value = "'9"
$ws.Range("G33").Value = "'9"
$ws.Range("G34").Value = "'9"
$ws.Range("G35").Value = "'9"
$ws.Range("G36").Value = "'9"
$ws.Range("G37").Value = "'9"
$ws.Range("G38").Value = "'9"
$ws.Range("D39").Value = "'0.01565"
$ws.Range("E39").Value = "'-1.89%"
$ws.Range("G39").Value = "'9"
$ws.Range("D40").Value = "'0.04444"
$ws.Range("E40").Value = "'0.74%"
$ws.Range("G40").Value = "'9"
$ws.Range("D41").Value = "'0.007402"
$ws.Range("E41").Value = "'1.77%"
$ws.Range("G41").Value = "'9"
$ws.Range("D42").Value = "'0.009215"
$ws.Range("E42").Value = "'-7.21%"
$ws.Range("G42").Value = "'9"
$ws.Range("E43").Value = "'1.54%"
$ws.Range("G43").Value = "'9"
$ws.Range("D44").Value = "'0.002222"
$ws.Range("E44").Value = "'7.15%"
$ws.Range("G44").Value = "'9"
$ws.Range("D45").Value = "'0.009089"
$ws.Range("E45").Value = "'-4.53%"
$ws.Range("G45").Value = "'9"
$ws.Range("D46").Value = "'0.00006114"
$ws.Range("E46").Value = "'2.66%"
$ws.Range("G46").Value = "'9"
$ws.Range("G47").Value = "'9"
$ws.Range("D48").Value = "'2.321"
$ws.Range("E48").Value = "'3.35%"
$ws.Range("G48").Value = "'9"
$ws.Range("D49").Value = "'0.002002"
$ws.Range("G49").Value = "'9"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("G50").Value = "'9"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("G51").Value = "'9"
